$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.382793479608949
$ws.Cells.Item(2, 3).Value = 0.295668839537484
$ws.Cells.Item(2, 4).Value = 0.0789221934265214
$ws.Cells.Item(2, 5).Value = 0.1015250414218549
$ws.Cells.Item(2, 7).Value = 0.4916681496755189
$ws.Cells.Item(2, 8).Value = 0.6313646754332041
$ws.Cells.Item(2, 12).Value = 0.2156223791713927
$ws.Cells.Item(2, 15).Value = 2.208208847287381
$ws.Cells.Item(3, 2).Value = 1.24211875765144
$ws.Cells.Item(3, 3).Value = 0.2811478313233522
$ws.Cells.Item(3, 4).Value = 0.07153395446266586
$ws.Cells.Item(3, 5).Value = 0.1028782273118711
$ws.Cells.Item(3, 7).Value = 0.4935663384005693
$ws.Cells.Item(3, 8).Value = 0.6373757017971755
$ws.Cells.Item(3, 12).Value = 0.2051410966078464
$ws.Cells.Item(3, 15).Value = 2.22447858884324
$ws.Cells.Item(4, 2).Value = 1.155679486509314
$ws.Cells.Item(4, 3).Value = 0.2722261428207275
$ws.Cells.Item(4, 4).Value = 0.0670323645140769
$ws.Cells.Item(4, 5).Value = 0.1037769102375972
$ws.Cells.Item(4, 7).Value = 0.4952528422650815
$ws.Cells.Item(4, 8).Value = 0.6414820881226717
$ws.Cells.Item(4, 12).Value = 0.1987996503366389
$ws.Cells.Item(4, 15).Value = 2.23643099991159
$ws.Cells.Item(5, 2).Value = 1.120440608262754
$ws.Cells.Item(5, 3).Value = 0.2685893578652951
$ws.Cells.Item(5, 4).Value = 0.06520666886756032
$ws.Cells.Item(5, 5).Value = 0.104160175165303
$ws.Cells.Item(5, 7).Value = 0.4960706752658623
$ws.Cells.Item(5, 8).Value = 0.6432598784631267
$ws.Cells.Item(5, 12).Value = 0.1962392017194219
$ws.Cells.Item(5, 15).Value = 2.241794092764351
$ws.Cells.Item(6, 2).Value = 1.114588423459338
$ws.Cells.Item(6, 3).Value = 0.2679854139736051
$ws.Cells.Item(6, 4).Value = 0.0649040408214745
$ws.Cells.Item(6, 5).Value = 0.104224845160239
$ws.Cells.Item(6, 7).Value = 0.4962143488585795
$ws.Cells.Item(6, 8).Value = 0.6435613820615913
$ws.Cells.Item(6, 12).Value = 0.1958154775099388
$ws.Cells.Item(6, 15).Value = 2.242714335403193
$ws.Cells.Item(7, 2).Value = 1.155204297358466
$ws.Cells.Item(7, 3).Value = 0.2721770999970659
$ws.Cells.Item(7, 4).Value = 0.06700770718185822
$ws.Cells.Item(7, 5).Value = 0.1037820100836928
$ws.Cells.Item(7, 7).Value = 0.4952633437670073
$ws.Cells.Item(7, 8).Value = 0.6415056414286369
$ws.Cells.Item(7, 12).Value = 0.1987650229932285
$ws.Cells.Item(7, 15).Value = 2.23650133633555
$ws.Cells.Item(8, 2).Value = 1.33430337196063
$ws.Cells.Item(8, 3).Value = 0.2906633694429104
$ws.Cells.Item(8, 4).Value = 0.0763674897045945
$ws.Cells.Item(8, 5).Value = 0.1019775377439149
$ws.Cells.Item(8, 7).Value = 0.492214242662854
$ws.Cells.Item(8, 8).Value = 0.6333509613289294
$ws.Cells.Item(8, 12).Value = 0.2119889452247605
$ws.Cells.Item(8, 15).Value = 2.213410594180402
$ws.Cells.Item(9, 2).Value = 1.684934013954944
$ws.Cells.Item(9, 3).Value = 0.3268574469618102
$ws.Cells.Item(9, 4).Value = 0.09500000021172639
$ws.Cells.Item(9, 5).Value = 0.0989775377044726
$ws.Cells.Item(9, 7).Value = 0.4903894979549221
$ws.Cells.Item(9, 8).Value = 0.6206618609687098
$ws.Cells.Item(9, 12).Value = 0.2386657811587725
$ws.Cells.Item(9, 15).Value = 2.18375678293782
$ws.Cells.Item(10, 2).Value = 1.942118810136094
$ws.Cells.Item(10, 3).Value = 0.3534009811071712
$ws.Cells.Item(10, 4).Value = 0.10886260616968
$ws.Cells.Item(10, 5).Value = 0.09710234318609245
$ws.Cells.Item(10, 7).Value = 0.4916108937332098
$ws.Cells.Item(10, 8).Value = 0.6133590231373915
$ws.Cells.Item(10, 12).Value = 0.2587190426737322
$ws.Cells.Item(10, 15).Value = 2.171573561198727
$ws.Cells.Item(11, 2).Value = 2.059014289801723
$ws.Cells.Item(11, 3).Value = 0.3654632774199058
$ws.Cells.Item(11, 4).Value = 0.1152075954398839
$ws.Cells.Item(11, 5).Value = 0.09632081285647587
$ws.Cells.Item(11, 7).Value = 0.4927291912584479
$ws.Cells.Item(11, 8).Value = 0.6104768371112925
$ws.Cells.Item(11, 12).Value = 0.267940437037069
$ws.Cells.Item(11, 15).Value = 2.168133006987318
$ws.Cells.Item(12, 2).Value = 2.103263703402547
$ws.Cells.Item(12, 3).Value = 0.3700288918943215
$ws.Cells.Item(12, 4).Value = 0.117615899830156
$ws.Cells.Item(12, 5).Value = 0.09603516151853242
$ws.Cells.Item(12, 7).Value = 0.493234066308986
$ws.Cells.Item(12, 8).Value = 0.6094488061665118
$ws.Cells.Item(12, 12).Value = 0.27144654366829
$ws.Cells.Item(12, 15).Value = 2.167133675383241
$ws.Cells.Item(13, 2).Value = 2.093734561696465
$ws.Cells.Item(13, 3).Value = 0.3690457042449395
$ws.Cells.Item(13, 4).Value = 0.1170969793728034
$ws.Cells.Item(13, 5).Value = 0.09609622343628921
$ws.Cells.Item(13, 7).Value = 0.4931217044020997
$ws.Cells.Item(13, 8).Value = 0.6096673896274041
$ws.Cells.Item(13, 12).Value = 0.2706908121696614
$ws.Cells.Item(13, 15).Value = 2.167335378765102
$ws.Cells.Item(14, 2).Value = 2.062655055700873
$ws.Cells.Item(14, 3).Value = 0.3658389374774913
$ws.Cells.Item(14, 4).Value = 0.1154056158969468
$ws.Cells.Item(14, 5).Value = 0.09629710569191197
$ws.Cells.Item(14, 7).Value = 0.4927690934690645
$ws.Cells.Item(14, 8).Value = 0.6103909891964605
$ws.Cells.Item(14, 12).Value = 0.2682286031169525
$ws.Cells.Item(14, 15).Value = 2.168044700627178
$ws.Cells.Item(15, 2).Value = 2.043615767228346
$ws.Cells.Item(15, 3).Value = 0.3638744153612947
$ws.Cells.Item(15, 4).Value = 0.1143703358806363
$ws.Cells.Item(15, 5).Value = 0.09642149331905792
$ws.Cells.Item(15, 7).Value = 0.4925637242157848
$ws.Cells.Item(15, 8).Value = 0.6108424740234426
$ws.Cells.Item(15, 12).Value = 0.2667222712605053
$ws.Cells.Item(15, 15).Value = 2.168518748506131
$ws.Cells.Item(16, 2).Value = 1.934477239371688
$ws.Cells.Item(16, 3).Value = 0.352612402940423
$ws.Cells.Item(16, 4).Value = 0.1084487284257705
$ws.Cells.Item(16, 5).Value = 0.09715485817049618
$ws.Cells.Item(16, 7).Value = 0.4915491746146614
$ws.Cells.Item(16, 8).Value = 0.6135562441722016
$ws.Cells.Item(16, 12).Value = 0.258118389388315
$ws.Cells.Item(16, 15).Value = 2.171840812795864
$ws.Cells.Item(17, 2).Value = 1.867497380193186
$ws.Cells.Item(17, 3).Value = 0.3457000953984561
$ws.Cells.Item(17, 4).Value = 0.1048259657600568
$ws.Cells.Item(17, 5).Value = 0.0976230785327008
$ws.Cells.Item(17, 7).Value = 0.4910712443169132
$ws.Cells.Item(17, 8).Value = 0.6153338198696048
$ws.Cells.Item(17, 12).Value = 0.2528655017754176
$ws.Cells.Item(17, 15).Value = 2.174418013852488
$ws.Cells.Item(18, 2).Value = 1.828963126389624
$ws.Cells.Item(18, 3).Value = 0.3417231620534835
$ws.Cells.Item(18, 4).Value = 0.1027459017306143
$ws.Cells.Item(18, 5).Value = 0.09789911600963741
$ws.Cells.Item(18, 7).Value = 0.4908492843904639
$ws.Cells.Item(18, 8).Value = 0.6163976330906422
$ws.Cells.Item(18, 12).Value = 0.2498535090063285
$ws.Cells.Item(18, 15).Value = 2.176098091004235
$ws.Cells.Item(19, 2).Value = 1.815914572854183
$ws.Cells.Item(19, 3).Value = 0.3403764518639605
$ws.Cells.Item(19, 4).Value = 0.1020422548543678
$ws.Cells.Item(19, 5).Value = 0.09799373290102942
$ws.Cells.Item(19, 7).Value = 0.4907832089689492
$ws.Cells.Item(19, 8).Value = 0.6167649282236738
$ws.Cells.Item(19, 12).Value = 0.2488353047891252
$ws.Cells.Item(19, 15).Value = 2.176700857995058
$ws.Cells.Item(20, 2).Value = 1.874628469120182
$ws.Cells.Item(20, 3).Value = 0.3464360443952899
$ws.Cells.Item(20, 4).Value = 0.1052112368443261
$ws.Cells.Item(20, 5).Value = 0.09757253908894192
$ws.Cells.Item(20, 7).Value = 0.4911166389127715
$ws.Cells.Item(20, 8).Value = 0.6151403083364073
$ws.Cells.Item(20, 12).Value = 0.2534237156751402
$ws.Cells.Item(20, 15).Value = 2.174123191571596
$ws.Cells.Item(21, 2).Value = 2.071784320870279
$ws.Cells.Item(21, 3).Value = 0.3667809025425299
$ws.Cells.Item(21, 4).Value = 0.1159022584150904
$ws.Cells.Item(21, 5).Value = 0.09623782212140775
$ws.Cells.Item(21, 7).Value = 0.4928704507979091
$ws.Cells.Item(21, 8).Value = 0.6101767290189173
$ws.Cells.Item(21, 12).Value = 0.268951429813356
$ws.Cells.Item(21, 15).Value = 2.16782810737493
$ws.Cells.Item(22, 2).Value = 2.20054070382696
$ws.Cells.Item(22, 3).Value = 0.3800649685867654
$ws.Cells.Item(22, 4).Value = 0.1229220744547348
$ws.Cells.Item(22, 5).Value = 0.0954255344724686
$ws.Cells.Item(22, 7).Value = 0.4944913548641239
$ws.Cells.Item(22, 8).Value = 0.6073022914234087
$ws.Cells.Item(22, 12).Value = 0.2791822145934049
$ws.Cells.Item(22, 15).Value = 2.165483727923942
$ws.Cells.Item(23, 2).Value = 2.131830506069264
$ws.Cells.Item(23, 3).Value = 0.3729762582917147
$ws.Cells.Item(23, 4).Value = 0.1191724782214436
$ws.Cells.Item(23, 5).Value = 0.09585357029264685
$ws.Cells.Item(23, 7).Value = 0.4935826534276799
$ws.Cells.Item(23, 8).Value = 0.6088025757285607
$ws.Cells.Item(23, 12).Value = 0.2737143275261928
$ws.Cells.Item(23, 15).Value = 2.166572594865102
$ws.Cells.Item(24, 2).Value = 1.871404586643905
$ws.Cells.Item(24, 3).Value = 0.3461033310381652
$ws.Cells.Item(24, 4).Value = 0.1050370473469968
$ws.Cells.Item(24, 5).Value = 0.09759536663123924
$ws.Cells.Item(24, 7).Value = 0.491095951556602
$ws.Cells.Item(24, 8).Value = 0.6152276645120907
$ws.Cells.Item(24, 12).Value = 0.2531713223659295
$ws.Cells.Item(24, 15).Value = 2.174255862781649
$ws.Cells.Item(25, 2).Value = 1.590148647062108
$ws.Cells.Item(25, 3).Value = 0.3170735254962267
$ws.Cells.Item(25, 4).Value = 0.08992925391349615
$ws.Cells.Item(25, 5).Value = 0.09973139977718937
$ws.Cells.Item(25, 7).Value = 0.4904353303068518
$ws.Cells.Item(25, 8).Value = 0.623740378746902
$ws.Cells.Item(25, 12).Value = 0.2313692862832823
$ws.Cells.Item(25, 15).Value = 2.167335378765102
